$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) values, columns A..Q
$headers = @(
    "button_dynamicObject_internalLabel",
    "button_dynamicObject_internalLabel_1",
    "button_dynamicObject_internalLabel_2",
    "button_dynamicObject_internalLabel_3",
    "button_dynamicObject_internalLabel_4",
    "button_dynamicObject_internalLabel_5",
    "button_dynamicObject_internalLabel_6",
    "button_dynamicObject_internalLabel_7",
    "button_dynamicObject_internalLabel_8",
    "button_dynamicObject_nthChild",
    "button_dynamicObject_nthChild_1",
    "button_dynamicObject_nthChild_2",
    "button_dynamicObject_nthChild_3",
    "button_dynamicObject_nthChild_4",
    "button_dynamicObject_nthChild_5",
    "input_object2",
    "link_moreOptions_internalRoleLinkName"
)

# Data row (row 2) text values for columns A..I and Q (P stays empty)
$textValues = @{
    "A2" = "Show slide 2 of"
    "B2" = "Show slide 3 of"
    "C2" = "Show slide 4 of"
    "D2" = "Show slide 5 of"
    "E2" = "Show slide 6 of"
    "F2" = "Show slide 7 of"
    "G2" = "Show slide 8 of"
    "H2" = "Show slide 9 of"
    "I2" = "Show slide 10 of"
    "Q2" = "+3 more"
}

# Data row (row 2) numeric-looking values that must stay as TEXT, columns J..O
$numericTextValues = @{
    "J2" = "3"
    "K2" = "3"
    "L2" = "3"
    "M2" = "4"
    "N2" = "4"
    "O2" = "3"
}

# Column widths, index 1..17 matching A..Q
$colWidths = @(36, 38, 38, 38, 38, 38, 38, 38, 38, 31, 33, 33, 33, 33, 33, 15, 39)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $headers[$i]
    $ws.Columns.Item($col).ColumnWidth = $colWidths[$i]
}

# Copy A1's formatting (the "Pandas" style) across the rest of row 1
$ws.Range("A1").Copy() | Out-Null
$ws.Range("B1:Q1").PasteSpecial(-4122) | Out-Null

foreach ($addr in $textValues.Keys) {
    $ws.Range($addr).Value = $textValues[$addr]
}

# Force these cells to hold text ("3"/"4") rather than being coerced to numbers
$ws.Range("J2:O2").NumberFormat = "@"
foreach ($addr in $numericTextValues.Keys) {
    $ws.Range($addr).Value = $numericTextValues[$addr]
}
